# Apply the change described in the diff:
#  - Swap the data in rows 2 and 3 (A2:G2 <-> A3:G3)
#  - Update the active selection to C15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = $ws.Range("A2:G2")
$row3 = $ws.Range("A3:G3")
$temp = $ws.Range("A20:G20")

# Move (not copy-by-value) the rows through a scratch range so the original
# cell typing/formatting travels with the data instead of being re-entered
# as a freshly-typed value (which would otherwise coerce text-like station
# numbers into numbers and disturb the shared-string table).
$row2.Cut($temp)
$row3.Cut($row2)
$temp.Cut($row3)

# Update the selected cell to match the diff.
$ws.Range("C15").Select()
